# Updated buckling safety factor
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Pin diameter (I20) is now derived from a 15/16" dimension converted to mm
# instead of being a hard-coded 5 mm value.
$ws.Range("I20").Formula = '=15/16*25.4'

# Euler buckling load formula now divides by the factor of safety in C20
# as well (previously only used C17). H21 becomes the new "master" cell of
# the shared formula group spanning H21:L21; G21 is written separately
# since it falls out of that shared range.
$ws.Range("G21").Formula = '=-G14*(G19/1000)^2/(PI()^2*$C$17/$C$20)'
$ws.Range("H21:L21").Formula = '=-H14*(H19/1000)^2/(PI()^2*$C$17/$C$20)'

# New check cell: reciprocal of the margin (in multiples of the 25.4 mm
# inch conversion) between the governing diameter (I25) and the pin
# diameter (I20).
$ws.Range("I34").Formula = '=((I25-I20)/25.4)^-1'

$ws.Range("J26").Select()
